# Added DP solution for SSS.
# Insert a new row for "17.0 Subset sum (SSS)" right before the existing
# "17.1 Ties in a presidential election" row (row 234), pushing the rest of
# section 17 down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 234 (shifts 234..246 down to 235..247).
$ws.Rows.Item(234).Insert()

# Fill in the new row's content.
$ws.Cells.Item(234, 2).Value = "17.0 Subset sum (SSS)"
$ws.Cells.Item(234, 3).Value = "****"
$ws.Cells.Item(234, 4).Value = "C1"

# Match formatting used by the rest of the section-17 rows: columns B/C use
# the plain bordered style, column D uses the "solved" (green) style, and
# column A picks up the section's shaded/no-border look.
$ws.Cells.Item(235, 2).Copy()
$ws.Cells.Item(234, 2).PasteSpecial(-4122)

$ws.Cells.Item(235, 3).Copy()
$ws.Cells.Item(234, 3).PasteSpecial(-4122)

$ws.Cells.Item(5, 4).Copy()
$ws.Cells.Item(234, 4).PasteSpecial(-4122)

$ws.Cells.Item(58, 2).Copy()
$ws.Cells.Item(234, 1).PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Restore the selection/scroll focus to where the new row now lives.
[void]$ws.Range("E234").Select()
